$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the extra last row (row 7) entirely so the used range shrinks to A1:D6
$ws.Rows(7).Delete()

# Clear out the Title/Location/Description columns for rows 3-6, keeping only
# the JobId values in column A.
$ws.Range("B3:D6").ClearContents()

# Row 6's JobId was mis-entered; correct it to 4.
$ws.Range("A6").Value = 4

$ws.Range("B3").Select()
